# "Added data for ProvarCache" -- Provar's RMA-number cache advances to a
# freshly generated pool (RMA-YZ8C-...) and the "RMA Details Maintenance
# Grid" sheet is repointed at that new pool (previously it referenced the
# RMA-VAXX-... pool).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 (line 1)
$ws.Range("E2").Value = "RMA-YZ8C-001"
$ws.Range("F2").Value = "RMA-YZ8C-1-1"
$ws.Range("J2").Value = "a7s5f000000xKM5AAM"

# Row 3 (line 2)
$ws.Range("E3").Value = "RMA-YZ8C-002"
$ws.Range("F3").Value = "RMA-YZ8C-1-2"
$ws.Range("J3").Value = "a7s5f000000xKM6AAM"

# Row 4 (line 3)
$ws.Range("E4").Value = "RMA-YZ8C-003"
$ws.Range("F4").Value = "RMA-YZ8C-1-3"
$ws.Range("J4").Value = "a7s5f000000xKM7AAM"
